{"js": "// Update the date line and the twenty-five \"two-digit \u00f7 one-digit\"\n// division problems to the next day's worksheet values.\nconst replacements = [\n  [\"2024-02-19 Monday\", \"2024-02-20 Tuesday\"],\n  [\"75\u00f72=\", \"27\u00f73=\"],\n  [\"71\u00f72=\", \"72\u00f79=\"],\n  [\"71\u00f78=\", \"65\u00f78=\"],\n  [\"98\u00f78=\", \"61\u00f75=\"],\n  [\"68\u00f74=\", \"14\u00f73=\"],\n  [\"56\u00f77=\", \"43\u00f79=\"],\n  [\"62\u00f73=\", \"99\u00f78=\"],\n  [\"41\u00f79=\", \"45\u00f73=\"],\n  [\"97\u00f75=\", \"75\u00f75=\"],\n  [\"47\u00f73=\", \"90\u00f75=\"],\n  [\"23\u00f72=\", \"54\u00f73=\"],\n  [\"71\u00f77=\", \"25\u00f74=\"],\n  [\"46\u00f78=\", \"56\u00f76=\"],\n  [\"57\u00f73=\", \"74\u00f78=\"],\n  [\"18\u00f74=\", \"53\u00f77=\"],\n  [\"76\u00f77=\", \"27\u00f79=\"],\n  [\"84\u00f76=\", \"92\u00f78=\"],\n  [\"98\u00f72=\", \"46\u00f75=\"],\n  [\"80\u00f77=\", \"93\u00f74=\"],\n  [\"73\u00f74=\", \"89\u00f75=\"],\n  [\"54\u00f72=\", \"80\u00f74=\"],\n  [\"21\u00f73=\", \"67\u00f72=\"],\n  [\"24\u00f79=\", \"36\u00f77=\"],\n  [\"44\u00f78=\", \"18\u00f78=\"],\n  [\"92\u00f79=\", \"14\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five \"two-digit \u00f7 one-digit\"\n# division problems to the next day's worksheet values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  , @(\"2024-02-19 Monday\", \"2024-02-20 Tuesday\")\n  , @(\"75\u00f72=\", \"27\u00f73=\")\n  , @(\"71\u00f72=\", \"72\u00f79=\")\n  , @(\"71\u00f78=\", \"65\u00f78=\")\n  , @(\"98\u00f78=\", \"61\u00f75=\")\n  , @(\"68\u00f74=\", \"14\u00f73=\")\n  , @(\"56\u00f77=\", \"43\u00f79=\")\n  , @(\"62\u00f73=\", \"99\u00f78=\")\n  , @(\"41\u00f79=\", \"45\u00f73=\")\n  , @(\"97\u00f75=\", \"75\u00f75=\")\n  , @(\"47\u00f73=\", \"90\u00f75=\")\n  , @(\"23\u00f72=\", \"54\u00f73=\")\n  , @(\"71\u00f77=\", \"25\u00f74=\")\n  , @(\"46\u00f78=\", \"56\u00f76=\")\n  , @(\"57\u00f73=\", \"74\u00f78=\")\n  , @(\"18\u00f74=\", \"53\u00f77=\")\n  , @(\"76\u00f77=\", \"27\u00f79=\")\n  , @(\"84\u00f76=\", \"92\u00f78=\")\n  , @(\"98\u00f72=\", \"46\u00f75=\")\n  , @(\"80\u00f77=\", \"93\u00f74=\")\n  , @(\"73\u00f74=\", \"89\u00f75=\")\n  , @(\"54\u00f72=\", \"80\u00f74=\")\n  , @(\"21\u00f73=\", \"67\u00f72=\")\n  , @(\"24\u00f79=\", \"36\u00f77=\")\n  , @(\"44\u00f78=\", \"18\u00f78=\")\n  , @(\"92\u00f79=\", \"14\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
